$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / description text for the "International" variant ---
$ws.Range("A1").Value = "2016 Internatioal Grant Dollars By Agency Table"
$ws.Range("A3").Value = "This table displays the total of international grant dollars each HHS agency awarded in FY 2016, provided as a text alternative to the interactive chart on the Awards page of this website."
$ws.Range("A5").Value = "HHS Total International Award Amounts Description"
$ws.Range("A7").Value = "Number of Iinternational grant dollars HHS awarded in FY 2016 by agency."

# --- Rename the "Organization" column header to "Agency" ---
$ws.Range("A9").Value = "Agency"

# --- Sort the agency/amount table (A10:B16) alphabetically by agency name ---
$dataRange = $ws.Range("A10:B16")
$keyRange = $ws.Range("A10:A16")
$dataRange.Sort($keyRange)

# --- Add a thin box border around the header row and the data table ---
$tableRange = $ws.Range("A9:B16")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# --- Move the active selection (matches the saved cursor position) ---
$ws.Range("I1").Select() | Out-Null
